$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 540
$ws.Range("E2").Value = 63
$ws.Range("F2").Value = 63
$ws.Range("G2").Value = 77
$ws.Range("H2").Value = 61
$ws.Range("I2").Value = 61
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = 1131
$ws.Range("L2").Value = 184
$ws.Range("M2").Value = 947
$ws.Range("N2").Value = 947
$ws.Range("O2").Value = ""
$ws.Range("P2").Value = 58
$ws.Range("Q2").Value = 111
$ws.Range("R2").Value = -121
$ws.Range("S2").Value = -27
$ws.Range("T2").Value = 153
$ws.Range("U2").Value = -42
$ws.Range("V2").Value = 10
$ws.Range("W2").Value = 11.59
$ws.Range("X2").Value = 11.32
$ws.Range("Y2").Value = 6.74
$ws.Range("Z2").Value = 5.7
$ws.Range("AA2").Value = 19.48
$ws.Range("AB2").Value = 1213.09
$ws.Range("AC2").Value = 526
$ws.Range("AD2").Value = 10.52
$ws.Range("AE2").Value = 8161
$ws.Range("AF2").Value = 0.68
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 2.71
$ws.Range("AI2").Value = 28.49
$ws.Range("AJ2").Value = 11600000

# Row 3
$ws.Range("D3").Value = 509
$ws.Range("E3").Value = 51
$ws.Range("F3").Value = 51
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = 54
$ws.Range("I3").Value = 54
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = 1179
$ws.Range("L3").Value = 202
$ws.Range("M3").Value = 977
$ws.Range("N3").Value = 977
$ws.Range("O3").Value = ""
$ws.Range("P3").Value = 58
$ws.Range("Q3").Value = 44
$ws.Range("R3").Value = -83
$ws.Range("S3").Value = 23
$ws.Range("T3").Value = 109
$ws.Range("U3").Value = -66
$ws.Range("V3").Value = 50
$ws.Range("W3").Value = 9.960000000000001
$ws.Range("X3").Value = 10.54
$ws.Range("Y3").Value = 5.58
$ws.Range("Z3").Value = 4.64
$ws.Range("AA3").Value = 20.69
$ws.Range("AB3").Value = 1269.35
$ws.Range("AC3").Value = 463
$ws.Range("AD3").Value = 13.9
$ws.Range("AE3").Value = 8424
$ws.Range("AF3").Value = 0.76
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 2.33
$ws.Range("AI3").Value = 32.43
$ws.Range("AJ3").Value = 11600000

# Row 4
$ws.Range("D4").Value = 513
$ws.Range("E4").Value = 38
$ws.Range("F4").Value = 38
$ws.Range("G4").Value = 48
$ws.Range("H4").Value = 38
$ws.Range("I4").Value = 38
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = 1131
$ws.Range("L4").Value = 175
$ws.Range("M4").Value = 956
$ws.Range("N4").Value = 956
$ws.Range("O4").Value = ""
$ws.Range("P4").Value = 58
$ws.Range("Q4").Value = 48
$ws.Range("R4").Value = -14
$ws.Range("S4").Value = -27
$ws.Range("T4").Value = 29
$ws.Range("U4").Value = 19
$ws.Range("V4").Value = 40
$ws.Range("W4").Value = 7.33
$ws.Range("X4").Value = 7.44
$ws.Range("Y4").Value = 3.95
$ws.Range("Z4").Value = 3.31
$ws.Range("AA4").Value = 18.25
$ws.Range("AB4").Value = 1300.22
$ws.Range("AC4").Value = 329
$ws.Range("AD4").Value = 15.94
$ws.Range("AE4").Value = 8244
$ws.Range("AF4").Value = 0.64
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 2.86
$ws.Range("AI4").Value = 45.55
$ws.Range("AJ4").Value = 11600000

# Row 5
$ws.Range("D5").Value = 489
$ws.Range("E5").Value = 18
$ws.Range("F5").Value = 18
$ws.Range("G5").Value = 20
$ws.Range("H5").Value = 16
$ws.Range("I5").Value = 16
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = 1107
$ws.Range("L5").Value = 132
$ws.Range("M5").Value = 975
$ws.Range("N5").Value = 975
$ws.Range("O5").Value = ""
$ws.Range("P5").Value = 58
$ws.Range("Q5").Value = 63
$ws.Range("R5").Value = -19
$ws.Range("S5").Value = -47
$ws.Range("T5").Value = 16
$ws.Range("U5").Value = 48
$ws.Range("V5").Value = 10
$ws.Range("W5").Value = 3.68
$ws.Range("X5").Value = 3.35
$ws.Range("Y5").Value = 1.7
$ws.Range("Z5").Value = 1.47
$ws.Range("AA5").Value = 13.55
$ws.Range("AB5").Value = 1304.22
$ws.Range("AC5").Value = 141
$ws.Range("AD5").Value = 32.31
$ws.Range("AE5").Value = 8401
$ws.Range("AF5").Value = 0.54
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 3.29
$ws.Range("AI5").Value = 106.16
$ws.Range("AJ5").Value = 11600000

# Row 6
$ws.Range("D6").Value = 474
$ws.Range("E6").Value = 32
$ws.Range("F6").Value = 32
$ws.Range("G6").Value = 46
$ws.Range("H6").Value = 29
$ws.Range("I6").Value = 29
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = 1232
$ws.Range("L6").Value = 152
$ws.Range("M6").Value = 1080
$ws.Range("N6").Value = 1080
$ws.Range("O6").Value = ""
$ws.Range("P6").Value = 58
$ws.Range("Q6").Value = 49
$ws.Range("R6").Value = -47
$ws.Range("S6").Value = -21
$ws.Range("T6").Value = 12
$ws.Range("U6").Value = 36
$ws.Range("V6").Value = 10
$ws.Range("W6").Value = 6.66
$ws.Range("X6").Value = 6.13
$ws.Range("Y6").Value = 2.83
$ws.Range("Z6").Value = 2.48
$ws.Range("AA6").Value = 14.07
$ws.Range("AB6").Value = 1355.17
$ws.Range("AC6").Value = 251
$ws.Range("AD6").Value = 18.52
$ws.Range("AE6").Value = 9310
$ws.Range("AF6").Value = 0.5
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 3.23
$ws.Range("AI6").Value = 59.48
$ws.Range("AJ6").Value = 11600000

# Rows 7-9: clear D:AJ, keep A-C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
